$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimates")

# Mark stories in rows 56, 65, 66 as Completed (column E of Table1 -
# the "Completed" checkbox column). This flips Completed Points (F)
# and Completed Hours (G) for those rows via their existing formulas.
$ws.Range("E56").Value = $true
$ws.Range("E65").Value = $true
$ws.Range("E66").Value = $true

# Bring the view to where the author left it: frozen pane scrolled so
# row 71 is the first visible row below the header, with E67 selected.
[void]$ws.Activate()
[void]$ws.Range("E67").Select()
$excel.ActiveWindow.ScrollRow = 71
